# Update the "想去人数" (want-to-go count) figures in column F for rows 2-5
# on both the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2" = 593
    "F3" = 3704
    "F4" = 103
    "F5" = 709
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
